$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, shifting existing rows 5-20 down to 6-21
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new entry
# (use single-quoted strings so the literal "$" is not treated as PowerShell variable interpolation)
$ws.Cells.Item(5, 1).Value = 'm4\~$Boss.xlsx'
$ws.Cells.Item(5, 2).Value = '~$Boss'
